$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Append the new bullet paragraphs describing the "Stellar et al., 2012"
#    operationalization coding entry (issues + resolution), at the end of
#    the document's list.
# ---------------------------------------------------------------------------

$quoteOpen  = [char]0x201C
$quoteClose = [char]0x201D

$textStellar    = "Stellar et al., 2012"
$textIssues     = "Issues:"
$textReasonA    = "Reason for operationalization in operationalization quote was not in reason quote (which was " + $quoteOpen + "None given" + $quoteClose + ") "
$textReasonB    = "(for both S2 and S3)"
$textGivenReason = "Given reason quote was not a reason for the operationalization (for both S2 and S3)"
$textResolution = "Resolution:"
$textCopied     = "Copied operationalization quote into reason cell (for both S2 and S3)"
$textDeleted    = "Deleted old reason quote (for both S2 and S3)"

$paraCountBefore = $d.Paragraphs.Count

$newTexts = @(
  $textStellar,
  $textIssues,
  $textReasonA,
  $textGivenReason,
  $textResolution,
  $textCopied,
  $textDeleted
)

$endOfDoc = $d.Content.End
$insertionRange = $d.Range($endOfDoc, $endOfDoc)
$joinedText = ($newTexts -join "`r")
$insertionRange.InsertAfter("`r" + $joinedText)

# New paragraphs are now at indices paraCountBefore+1 .. paraCountBefore+7,
# all cloned (pStyle/rPr/ilvl=3/numId=1) from the paragraph that used to be
# last. Fix the outline (list) level for the ones that are not level 4
# (ilvl=3): "Stellar et al., 2012" -> ilvl 1, "Issues:" -> ilvl 2,
# "Resolution:" -> ilvl 2 (ListLevelNumber is 1-based, i.e. ilvl + 1).

$pStellar = $d.Paragraphs($paraCountBefore + 1)
$pStellar.Range.ListFormat.ListLevelNumber = 2

$pIssues = $d.Paragraphs($paraCountBefore + 2)
$pIssues.Range.ListFormat.ListLevelNumber = 3

$pResolution = $d.Paragraphs($paraCountBefore + 5)
$pResolution.Range.ListFormat.ListLevelNumber = 3

# ---------------------------------------------------------------------------
# 2. The "Reason for operationalization ... (which was "None given")" bullet
#    is made of two separate runs (identical formatting) in the target, so
#    split off "(for both S2 and S3)" into its own run by inserting it
#    after the first chunk, then nudging a toggled/reverted character
#    property across just that new span to force a run break without
#    leaving any residual direct formatting behind.
# ---------------------------------------------------------------------------

$pReason = $d.Paragraphs($paraCountBefore + 3)
$reasonParaStart = $pReason.Range.Start
$splitPos = $reasonParaStart + $textReasonA.Length

$insB = $d.Range($splitPos, $splitPos)
$insB.InsertAfter($textReasonB)

$newRunEnd = $splitPos + $textReasonB.Length
$newRunRange = $d.Range($splitPos, $newRunEnd)
$newRunRange.Bold = 1
$newRunRange.Bold = 0

# ---------------------------------------------------------------------------
# 3. Add the "ListLabel 82".."ListLabel 90" character styles that the
#    original authoring tool minted alongside the new list paragraphs.
# ---------------------------------------------------------------------------

function Add-ListLabelStyle($num, $hasLatin, $hasSize) {
  $style = $d.Styles.Add("ListLabel " + $num, 2)
  if ($hasLatin) {
    $style.Font.NameAscii = "Times New Roman"
    $style.Font.NameOther = "Times New Roman"
  }
  $style.Font.NameBi = "OpenSymbol"
  if ($hasSize) {
    $style.Font.Size = 12
  }
  $style.QuickStyle = $true
}

Add-ListLabelStyle 82 $true  $false
Add-ListLabelStyle 83 $true  $true
Add-ListLabelStyle 84 $true  $true
Add-ListLabelStyle 85 $true  $true
Add-ListLabelStyle 86 $false $false
Add-ListLabelStyle 87 $false $false
Add-ListLabelStyle 88 $false $false
Add-ListLabelStyle 89 $false $false
Add-ListLabelStyle 90 $false $false
